$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell C1 (13-01-2023), formatted the same as B1 (bold/bordered/centered)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats

# Swap the "total" / "Alpha planeam equil" labels between row 2 and row 4
$ws.Range("A2").Value = "Alpha planeam equil"
$ws.Range("A4").Value = "total"

# New column C values (plain, unstyled, matching column B's numeric cells)
$ws.Range("C2").Value = 461.21
$ws.Range("C3").Value = 461.21
$ws.Range("C4").Value = 461.21
